$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this data block (rows 679-680),
# pushing the existing rows 679.. down by two (old 679 -> 681, ... old 728 -> 730)
$ws.Rows("679:680").Insert()

# New row 679: new weekly record (same variety/quality/volume as the old row 679 had,
# i.e. Valencia / Primera / 400), with updated date and prices.
$ws.Range("A679").Value = 4
$ws.Range("B679").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C679").Value = "Los Lagos"
$ws.Range("D679").Value = 45013
$ws.Range("E679").Value = 10
$ws.Range("F679").Value = "Fruta"
$ws.Range("G679").Value = 100102
$ws.Range("H679").Value = "Cítricos"
$ws.Range("I679").Value = 100102005
$ws.Range("J679").Value = "Naranja"
$ws.Range("K679").Value = "Valencia"
$ws.Range("L679").Value = "Primera"
$ws.Range("M679").Value = 400
$ws.Range("N679").Value = 17000
$ws.Range("O679").Value = 18000
$ws.Range("P679").Value = 17500
$ws.Range("Q679").Value = "$/caja 15 kilos empedrada"
$ws.Range("R679").Value = "Región de O'Higgins"
$ws.Range("S679").Value = 1167
$ws.Range("T679").Value = 15

# New row 680: another new weekly record (Valencia / Segunda / 200).
$ws.Range("A680").Value = 4
$ws.Range("B680").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C680").Value = "Los Lagos"
$ws.Range("D680").Value = 45013
$ws.Range("E680").Value = 10
$ws.Range("F680").Value = "Fruta"
$ws.Range("G680").Value = 100102
$ws.Range("H680").Value = "Cítricos"
$ws.Range("I680").Value = 100102005
$ws.Range("J680").Value = "Naranja"
$ws.Range("K680").Value = "Valencia"
$ws.Range("L680").Value = "Segunda"
$ws.Range("M680").Value = 200
$ws.Range("N680").Value = 15000
$ws.Range("O680").Value = 15000
$ws.Range("P680").Value = 15000
$ws.Range("Q680").Value = "$/caja 15 kilos empedrada"
$ws.Range("R680").Value = "Región de O'Higgins"
$ws.Range("S680").Value = 1000
$ws.Range("T680").Value = 15
